# --------------------------------------------------------------------------
# GitHub Actions refresh: "Updated cryptos list on Mon Jul 22 13:52:58 UTC
# 2024 with GitHub Actions".
#
# Re-scrapes coinranking.com and rewrites each coin row's Price (col D) and
# Volume(1h) (col E) with the latest quote. LEO (rank 24) overtook Polygon
# (rank 23) in this run, so rows 25-26 also swap their Coin/Link/Price/Volume
# contents to reflect the new ranking order.
# --------------------------------------------------------------------------
#
# NOTE: the Price column stores numeric-looking quotes ("594.68", "1.00", …)
# as literal text (that's how the scraper originally wrote them), so values
# that Excel would otherwise auto-convert to a Number are entered with a
# leading apostrophe (e.g. '594.68) -- exactly like typing '594.68 into a
# cell by hand -- to keep them text. Values that already contain two dots
# (e.g. "67.239.49") aren't valid numbers, so no apostrophe is needed there.
# --------------------------------------------------------------------------

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2: Bitcoin
$ws.Range("D2").Value = '67.239.49'
$ws.Range("E2").Value = '  +0.78%  '

# Row 3: Ethereum
$ws.Range("D3").Value = '3.472.78'
$ws.Range("E3").Value = '  -0.51%  '

# Row 4: TetherUSD
$ws.Range("E4").Value = '  +0.04%  '

# Row 5: BNB
$ws.Range("D5").Value = '''594.68'
$ws.Range("E5").Value = '  +0.19%  '

# Row 6: Solana
$ws.Range("D6").Value = '''179.27'
$ws.Range("E6").Value = '  +4.24%  '

# Row 7: XRP
$ws.Range("D7").Value = '''0.608'
$ws.Range("E7").Value = '  +5.14%  '

# Row 9: LidoStakedEther
$ws.Range("D9").Value = '3.472.47'
$ws.Range("E9").Value = '  -0.38%  '

# Row 10: Dogecoin
$ws.Range("D10").Value = '''0.138'
$ws.Range("E10").Value = '  +5.19%  '

# Row 11: Toncoin
$ws.Range("D11").Value = '''7.00'
$ws.Range("E11").Value = '  -1.58%  '

# Row 12: Cardano
$ws.Range("D12").Value = '''0.433'
$ws.Range("E12").Value = '  +0.99%  '

# Row 13: WrappedliquidstakedEther2.0
$ws.Range("D13").Value = '4.088.38'
$ws.Range("E13").Value = '  -0.17%  '

# Row 14: Avalanche
$ws.Range("D14").Value = '''31.72'
$ws.Range("E14").Value = '  +8.82%  '

# Row 15: TRON
$ws.Range("D15").Value = '''0.134'
$ws.Range("E15").Value = '  +0.01%  '

# Row 16: WrappedBTC
$ws.Range("D16").Value = '67.292.44'
$ws.Range("E16").Value = '  +0.86%  '

# Row 17: ShibaInu
$ws.Range("D17").Value = '''0.0000177'
$ws.Range("E17").Value = '  -0.41%  '

# Row 18: WrappedEther
$ws.Range("D18").Value = '3.490.07'
$ws.Range("E18").Value = '  +0.04%  '

# Row 19: Polkadot
$ws.Range("D19").Value = '''6.28'
$ws.Range("E19").Value = '  +0.57%  '

# Row 20: Chainlink
$ws.Range("D20").Value = '''14.16'
$ws.Range("E20").Value = '  -0.66%  '

# Row 21: BitcoinCash
$ws.Range("D21").Value = '''388.44'
$ws.Range("E21").Value = '  -0.27%  '

# Row 22: Uniswap
$ws.Range("D22").Value = '''7.91'
$ws.Range("E22").Value = '  +0.26%  '

# Row 23: Litecoin
$ws.Range("D23").Value = '''72.71'
$ws.Range("E23").Value = '  -0.74%  '

# Row 24: Dai
$ws.Range("D24").Value = '''1.00'
$ws.Range("E24").Value = '  +0.12%  '

# Row 27: PEPE
$ws.Range("D27").Value = '''0.0000122'
$ws.Range("E27").Value = '  +1.62%  '

# Row 28: InternetComputer(DFINITY)
$ws.Range("D28").Value = '''10.28'
$ws.Range("E28").Value = '  +2.23%  '

# Row 29: Kaspa
$ws.Range("D29").Value = '''0.175'
$ws.Range("E29").Value = '  -2.34%  '

# Row 30: Binance-PegBSC-USD
$ws.Range("D30").Value = '''1.00'
$ws.Range("E30").Value = '  +0.42%  '

# Row 31: NEARProtocol
$ws.Range("D31").Value = '''6.14'
$ws.Range("E31").Value = '  +0.72%  '

# Row 32: Fetch.AI
$ws.Range("E32").Value = '  +0.43%  '

# Row 33: PancakeSwap
$ws.Range("D33").Value = '''2.06'
$ws.Range("E33").Value = '  +0.76%  '

# Row 34: EthereumClassic
$ws.Range("D34").Value = '''23.47'
$ws.Range("E34").Value = '  -0.35%  '

# Row 35: Aptos
$ws.Range("D35").Value = '''7.36'
$ws.Range("E35").Value = '  +0.42%  '

# Row 36: USDe
$ws.Range("E36").Value = '  +0.03%  '

# Row 37: ImmutableX
$ws.Range("D37").Value = '''1.60'
$ws.Range("E37").Value = '  +0.50%  '

# Row 38: Monero
$ws.Range("D38").Value = '''161.96'
$ws.Range("E38").Value = '  -1.19%  '

# Row 39: Mantle
$ws.Range("D39").Value = '''0.881'
$ws.Range("E39").Value = '  +1.11%  '

# Row 40: dogwifhat
$ws.Range("D40").Value = '''2.82'
$ws.Range("E40").Value = '  +11.23%  '

# Row 41: Stacks
$ws.Range("D41").Value = '''1.87'
$ws.Range("E41").Value = '  -1.21%  '

# Row 42: RenderToken
$ws.Range("D42").Value = '''6.82'
$ws.Range("E42").Value = '  +0.23%  '

# Row 43: Filecoin
$ws.Range("D43").Value = '''4.61'
$ws.Range("E43").Value = '  -0.21%  '

# Row 44: EnergySwap
$ws.Range("D44").Value = '''26.14'
$ws.Range("E44").Value = '  +1.06%  '

# Row 45: Maker
$ws.Range("D45").Value = '2.808.45'
$ws.Range("E45").Value = '  -0.35%  '

# Row 46: InjectiveProtocol
$ws.Range("D46").Value = '''26.61'
$ws.Range("E46").Value = '  -1.26%  '

# Row 47: Hedera
$ws.Range("D47").Value = '''0.0721'
$ws.Range("E47").Value = '  -0.66%  '

# Row 48: OKB
$ws.Range("D48").Value = '''41.23'
$ws.Range("E48").Value = '  -2.83%  '

# Row 49: VeChain
$ws.Range("D49").Value = '''0.0299'
$ws.Range("E49").Value = '  +0.21%  '

# Row 50: Bittensor
$ws.Range("D50").Value = '''331.53'
$ws.Range("E50").Value = '  -1.87%  '

# Row 51: ONDO
$ws.Range("D51").Value = '''1.05'
$ws.Range("E51").Value = '  -1.49%  '

# Rows 25-26: LEO climbed past Polygon this cycle, so the two coins swap rows
# (Coin name, Link, Price and Volume all move together) in addition to each
# picking up its own freshly scraped Price/Volume figures.
# Row 25 was Polygon, now LEO:
$ws.Range("B25").Value = 'LEO'
$ws.Range("C25").Value = 'https://coinranking.com/coin/mqtUpyBxu8O8+leo-leo'
$ws.Range("D25").Value = '''5.76'
$ws.Range("E25").Value = '  +1.45%  '

# Row 26 was LEO, now Polygon:
$ws.Range("B26").Value = 'Polygon'
$ws.Range("C26").Value = 'https://coinranking.com/coin/uW2tk-ILY0ii+polygon-matic'
$ws.Range("D26").Value = '''0.538'
$ws.Range("E26").Value = '  +1.14%  '
